$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: "Top 10 Cities" - replace "philadelphia" row with "fort worth"
# appended at the bottom, shifting rows 4-10 up by one (row 4 removed,
# rows 5-11 become rows 4-10, new row 11 = fort worth)
# ---------------------------------------------------------------
$wsCities = $wb.Worksheets.Item("Top 10 Cities")

$wsCities.Range("A4").Value = "austin"
$wsCities.Range("B4").Value = "Texas"
$wsCities.Range("C4").Value = 6277143000
$wsCities.Range("D4").Value = 6525.775028589251

$wsCities.Range("A5").Value = "phoenix"
$wsCities.Range("B5").Value = "Arizona"
$wsCities.Range("C5").Value = 4744660000
$wsCities.Range("D5").Value = 2950.31059762839

$wsCities.Range("A6").Value = "houston"
$wsCities.Range("B6").Value = "Texas"
$wsCities.Range("C6").Value = 4542769000
$wsCities.Range("D6").Value = 1973.767929050232

$wsCities.Range("A7").Value = "portland"
$wsCities.Range("B7").Value = "Oregon"
$wsCities.Range("C7").Value = 4381080492
$wsCities.Range("D7").Value = 6714.114387649076

$wsCities.Range("A8").Value = "boston"
$wsCities.Range("B8").Value = "Massachusetts"
$wsCities.Range("C8").Value = 3390433000
$wsCities.Range("D8").Value = 5018.165214199445

$wsCities.Range("A9").Value = "jacksonville"
$wsCities.Range("B9").Value = "Florida"
$wsCities.Range("C9").Value = 3244209000
$wsCities.Range("D9").Value = 3416.337761790029

$wsCities.Range("A10").Value = "san jose"
$wsCities.Range("B10").Value = "California"
$wsCities.Range("C10").Value = 3116478000
$wsCities.Range("D10").Value = 3075.812680550443

$wsCities.Range("A11").Value = "fort worth"
$wsCities.Range("B11").Value = "Texas"
$wsCities.Range("C11").Value = 2926444000
$wsCities.Range("D11").Value = 3184.596985429904

# ---------------------------------------------------------------
# Sheet: "Top 10 Cities PC" - last row philadelphia -> miami
# ---------------------------------------------------------------
$wsCitiesPC = $wb.Worksheets.Item("Top 10 Cities PC")

$wsCitiesPC.Range("A11").Value = "miami"
$wsCitiesPC.Range("B11").Value = "Florida"
$wsCitiesPC.Range("C11").Value = 4700.251071767739
$wsCitiesPC.Range("D11").Value = 2078733039

# ---------------------------------------------------------------
# Sheet: "Top 10 Schools Least Debt"
# ---------------------------------------------------------------
$wsSchoolsLeastDebt = $wb.Worksheets.Item("Top 10 Schools Least Debt")

$wsSchoolsLeastDebt.Range("B2").Value = "rutherford county schools"
$wsSchoolsLeastDebt.Range("C2").Value = "Tennessee"
$wsSchoolsLeastDebt.Range("D2").Value = -76790573
$wsSchoolsLeastDebt.Range("E2").Value = -1559.104480945323

$wsSchoolsLeastDebt.Range("D3").Value = -65971070
$wsSchoolsLeastDebt.Range("E3").Value = -885.6841554117552

$wsSchoolsLeastDebt.Range("B7").Value = "detroit public schools"
$wsSchoolsLeastDebt.Range("C7").Value = "Michigan"

$wsSchoolsLeastDebt.Range("B8").Value = "hawaii department of education"
$wsSchoolsLeastDebt.Range("C8").Value = "Hawaii"

$wsSchoolsLeastDebt.Range("B9").Value = "board of education of baltimore county"
$wsSchoolsLeastDebt.Range("C9").Value = "Maryland"
$wsSchoolsLeastDebt.Range("D9").Value = 1844637
$wsSchoolsLeastDebt.Range("E9").Value = 16.59801504463

$wsSchoolsLeastDebt.Range("B10").Value = "chesterfield county school board"
$wsSchoolsLeastDebt.Range("C10").Value = "Virginia"
$wsSchoolsLeastDebt.Range("D10").Value = 4167372
$wsSchoolsLeastDebt.Range("E10").Value = 66.73668027864521

$wsSchoolsLeastDebt.Range("B11").Value = "city and county of denver school district no. 1"
$wsSchoolsLeastDebt.Range("C11").Value = "Colorado"
$wsSchoolsLeastDebt.Range("D11").Value = 30587841
$wsSchoolsLeastDebt.Range("E11").Value = 344.027634375949
